# Fix -- routes + position for btns import and export
#
# The "Catégories" import template was missing a column telling the
# importer WHICH sheet/model a parent category applies to (Balles /
# Projecteur / Raquettes démoniaques all feed the "Matériel" import).
# Add that "Catégorie pour" column, and leave the workbook focused on
# the sheet/cell the user was last working in.

$wb = $excel.ActiveWorkbook

$wsCategories = $wb.Worksheets.Item("Catégories")

# New column C: "Catégorie pour" header + "Matériel" for every existing row.
$wsCategories.Range("C1").Value2 = "Catégorie pour"
for ($row = 2; $row -le 7; $row++) {
    $wsCategories.Cells.Item($row, 3).Value2 = "Matériel"
}

# Leave "Catégories" as the active sheet/tab, with the cursor parked just
# below the data (next free row) ready for the next import row.
$wsCategories.Activate() | Out-Null
$wsCategories.Range("C9").Select() | Out-Null
